# Refresh the run timestamps recorded in column Z ("timestamp") of the
# Log_Muestras sheet. These values are written once per pcSMOTE sample log
# entry at the moment the augmentation run executes, so re-running the
# pipeline (e.g. after adding the "Us Crime" dataset run) updates every
# existing row's timestamp to the new run time while leaving all other
# columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z3").Value = "2025-11-13T06:52:25.315704"
$ws.Range("Z4:Z9").Value = "2025-11-13T06:52:25.316710"
$ws.Range("Z10:Z11").Value = "2025-11-13T06:52:25.317708"
$ws.Range("Z12:Z49").Value = "2025-11-13T06:52:25.318211"
$ws.Range("Z50").Value = "2025-11-13T06:52:25.323814"
$ws.Range("Z51:Z52").Value = "2025-11-13T06:52:25.324836"
$ws.Range("Z53:Z57").Value = "2025-11-13T06:52:25.325813"
$ws.Range("Z58:Z60").Value = "2025-11-13T06:52:25.429847"
$ws.Range("Z61:Z70").Value = "2025-11-13T06:52:25.430849"
$ws.Range("Z71:Z73").Value = "2025-11-13T06:52:25.586989"
$ws.Range("Z74:Z79").Value = "2025-11-13T06:52:25.587988"
